$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark (currently sitting by itself in the
#    empty paragraph right before the final section break)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Re-insert the _GoBack bookmark right after "Prof. " in the title
#    paragraph; this splits the "Prof. Ing. Roberto Martínez Román" run into
#    "Prof. " and "Ing. Roberto Martínez Román"
$p = $d.Paragraphs(2)
$r = $p.Range
$r.Find.Execute("Prof. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. Drop "Ing. " from the (now separate) second run so it reads
#    "Roberto Martínez Román"
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.Find.Execute("Ing. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
